$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.9610983981693364
$ws.Range("D2").Value = 0.9614686212176559
$ws.Range("E2").Value = 0.9610983981693364
$ws.Range("F2").Value = 0.9602613857752854

$ws.Range("C3").Value = 0.9794050343249427
$ws.Range("D3").Value = 0.9795578020533958
$ws.Range("E3").Value = 0.9794050343249427
$ws.Range("F3").Value = 0.9790161131238534

$ws.Range("C4").Value = 0.9988558352402745
$ws.Range("D4").Value = 0.9988571491136828
$ws.Range("E4").Value = 0.9988558352402745
$ws.Range("F4").Value = 0.9987781526637115

$ws.Range("C5").Value = 0.9590007627765065
$ws.Range("D5").Value = 0.9590672423841434
$ws.Range("E5").Value = 0.9590007627765065
$ws.Range("F5").Value = 0.9589334073701533

$ws.Range("C6").Value = 0.982837528604119
$ws.Range("D6").Value = 0.9830387022797661
$ws.Range("E6").Value = 0.982837528604119
$ws.Range("F6").Value = 0.981974680149756

$ws.Range("C7").Value = 0.9893211289092296
$ws.Range("D7").Value = 0.9893029321775226
$ws.Range("E7").Value = 0.9893211289092296
$ws.Range("F7").Value = 0.9888442183535571

$ws.Range("C8").Value = 0.9937070938215103
$ws.Range("D8").Value = 0.9936993144018478
$ws.Range("E8").Value = 0.9937070938215103
$ws.Range("F8").Value = 0.9933860189153357

$ws.Range("C9").Value = 0.9975209763539283
$ws.Range("D9").Value = 0.9974982465586499
$ws.Range("E9").Value = 0.9975209763539283
$ws.Range("F9").Value = 0.9974386488620747

$ws.Range("C10").Value = 0.9950419527078566
$ws.Range("D10").Value = 0.9950334815486824
$ws.Range("E10").Value = 0.9950419527078566
$ws.Range("F10").Value = 0.9948618772884045

$ws.Range("C12").Value = 0.9900839054157132
$ws.Range("D12").Value = 0.9901883276349078
$ws.Range("E12").Value = 0.9900839054157132
$ws.Range("F12").Value = 0.9897218683393449

$ws.Range("C13").Value = 0.9879862700228833
$ws.Range("D13").Value = 0.9879331871168294
$ws.Range("E13").Value = 0.9879862700228833
$ws.Range("F13").Value = 0.987792648483144

$ws.Range("C14").Value = 0.9858886346300534
$ws.Range("D14").Value = 0.985890248708102
$ws.Range("E14").Value = 0.9858886346300534
$ws.Range("F14").Value = 0.9854588633226952

$ws.Range("C15").Value = 0.9973302822273074
$ws.Range("D15").Value = 0.9973375158345879
$ws.Range("E15").Value = 0.9973302822273074
$ws.Range("F15").Value = 0.997220854841857

$ws.Range("C16").Value = 0.9958047292143402
$ws.Range("D16").Value = 0.9957829141587861
$ws.Range("E16").Value = 0.9958047292143402
$ws.Range("F16").Value = 0.9956215277043468

$ws.Range("C17").Value = 0.9969488939740656
$ws.Range("D17").Value = 0.9969582928374492
$ws.Range("E17").Value = 0.9969488939740656
$ws.Range("F17").Value = 0.996740826209381

$ws.Range("C18").Value = 0.9956140350877193
$ws.Range("D18").Value = 0.9956338265715299
$ws.Range("E18").Value = 0.9956140350877193
$ws.Range("F18").Value = 0.9954598818556766

$ws.Range("C19").Value = 0.9918001525553013
$ws.Range("D19").Value = 0.9918703343706122
$ws.Range("E19").Value = 0.9918001525553013
$ws.Range("F19").Value = 0.9914527698556737

$ws.Range("C20").Value = 0.9729214340198322
$ws.Range("D20").Value = 0.9737448079746895
$ws.Range("E20").Value = 0.9729214340198322
$ws.Range("F20").Value = 0.9716400846631509

$ws.Range("C21").Value = 0.9864607170099161
$ws.Range("D21").Value = 0.9865962758749373
$ws.Range("E21").Value = 0.9864607170099161
$ws.Range("F21").Value = 0.985730747981579

$ws.Range("C22").Value = 0.988558352402746
$ws.Range("D22").Value = 0.9886318138113294
$ws.Range("E22").Value = 0.988558352402746
$ws.Range("F22").Value = 0.9878605666067279

$ws.Range("C23").Value = 0.9900839054157132
$ws.Range("D23").Value = 0.9900893476441833
$ws.Range("E23").Value = 0.9900839054157132
$ws.Range("F23").Value = 0.9896077043630818

$ws.Range("C24").Value = 0.9973302822273074
$ws.Range("D24").Value = 0.9973375298613616
$ws.Range("E24").Value = 0.9973302822273074
$ws.Range("F24").Value = 0.9972326922134557

$ws.Range("C25").Value = 0.9986651411136537
$ws.Range("D25").Value = 0.9986669345900845
$ws.Range("E25").Value = 0.9986651411136537
$ws.Range("F25").Value = 0.9986032963694753

$ws.Range("C26").Value = 0.9984744469870328
$ws.Range("D26").Value = 0.9984767939916682
$ws.Range("E26").Value = 0.9984744469870328
$ws.Range("F26").Value = 0.9984114694810127

$ws.Range("C27").Value = 0.9988558352402745
$ws.Range("D27").Value = 0.9988571491136828
$ws.Range("E27").Value = 0.9988558352402745
$ws.Range("F27").Value = 0.9987781526637115

$ws.Range("C28").Value = 0.9979023646071701
$ws.Range("D28").Value = 0.9979068070386914
$ws.Range("E28").Value = 0.9979023646071701
$ws.Range("F28").Value = 0.9977995395035014

$ws.Range("C29").Value = 0.9900839054157132
$ws.Range("D29").Value = 0.9901272713322234
$ws.Range("E29").Value = 0.9900839054157132
$ws.Range("F29").Value = 0.9895068038989967

$ws.Range("C30").Value = 0.9731121281464531
$ws.Range("D30").Value = 0.9733276997163469
$ws.Range("E30").Value = 0.9731121281464531
$ws.Range("F30").Value = 0.9728516613925116

$ws.Range("C31").Value = 0.9887490465293669
$ws.Range("D31").Value = 0.9887600503347537
$ws.Range("E31").Value = 0.9887490465293669
$ws.Range("F31").Value = 0.9884351757961887

$ws.Range("C32").Value = 0.9874141876430206
$ws.Range("D32").Value = 0.987258216629961
$ws.Range("E32").Value = 0.9874141876430206
$ws.Range("F32").Value = 0.9871746138310818

$ws.Range("C33").Value = 0.9986651411136537
$ws.Range("D33").Value = 0.9986669390763637
$ws.Range("E33").Value = 0.9986651411136537
$ws.Range("F33").Value = 0.998619333421091

$ws.Range("C34").Value = 0.9935163996948894
$ws.Range("D34").Value = 0.9934801213779205
$ws.Range("E34").Value = 0.9935163996948894
$ws.Range("F34").Value = 0.9934881668919573

$ws.Range("D35").Value = 0.9963536549221774
$ws.Range("F35").Value = 0.9962295767840115

$ws.Range("C36").Value = 0.9887490465293669
$ws.Range("D36").Value = 0.9888246208345635
$ws.Range("E36").Value = 0.9887490465293669
$ws.Range("F36").Value = 0.9881307176983142

$ws.Range("C37").Value = 0.9731121281464531
$ws.Range("D37").Value = 0.9733036207132412
$ws.Range("E37").Value = 0.9731121281464531
$ws.Range("F37").Value = 0.9725707615990664
